$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A46").Value = "http://purl.obolibrary.org/obo/GO_0008150 "
$ws.Range("B46").Value = "biological_process"
$ws.Range("C46").Value = "y"

$ws.Range("A47").Value = "http://purl.obolibrary.org/obo/IAO_0000416"
$ws.Range("A48").Value = "http://purl.obolibrary.org/obo/OBI_0001619"
$ws.Range("B48").Value = "specimen collection time measurement datum"
$ws.Range("B47").Value = "time measurement datum"
$ws.Range("C47").Value = "y"
$ws.Range("C48").Value = "y"

$ws.Range("A47").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
